{"js": "// Office.js (Word JavaScript API) script.\n// Applies the documented edit to READ_THIS_FIRST.docx:\n//  1. \"During the install...\" paragraph: drop \"Client only or \" and append a\n//     new sentence about creating a root account/password.\n//  2. Remove the old \"Step 2\" body paragraph (the MySQL-Client-setup\n//     hyperlink blurb) together with the whole \"Step 3\" / \"Step 4\" sections,\n//     because their replacement content already exists lower in the\n//     document (the paragraphs that used to live under \"Step 4\") - removing\n//     the now-obsolete paragraphs promotes that content up to sit right\n//     under \"Step 2\".\n//  3. Tweak the (promoted) \"Once this is opened, open a new SQL script\n//     file...\" paragraph: insert \"(not create!)\" and the \"located in the\n//     same directory...\" clause.\n//  4. Renumber the final remaining Heading2 (\"Step 5\") to \"Step 3\".\n\nconst body = context.document.body;\n\n// --- 1. \"During the install...\" paragraph -------------------------------\nconst clientOnly = body.search(\"Client only or \", { matchCase: true });\nclientOnly.load(\"items\");\nawait context.sync();\nif (clientOnly.items.length > 0) {\n  clientOnly.items[0].insertText(\"\", \"Replace\");\n}\n\nconst selectedPeriod = body.search(\"Full must be selected. \", { matchCase: true });\nselectedPeriod.load(\"items\");\nawait context.sync();\nif (selectedPeriod.items.length > 0) {\n  selectedPeriod.items[0].insertText(\n    \"Full must be selected. Make sure that in the setup process you create a root account and password for the account.\",\n    \"Replace\"\n  );\n}\n\n// --- 2. Drop the obsolete paragraphs/headings ----------------------------\n// Order in the (current) document right after the \"Step 2\" heading:\n//   \"Once MySQL has been installed, follow the direction at this website...\"\n//   Heading2 \"Step 3\"\n//   \"If one is not already present, create a new connection...\"\n//   Heading2 \"Step 4\"\n// All four are removed; the paragraphs that used to sit beneath \"Step 4\"\n// (root-connection / SQL script / schema / exit instructions) shift up to\n// immediately follow \"Step 2\".\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (\n    t.indexOf(\"Once MySQL has been installed, follow the direction at\") === 0 ||\n    t === \"Step 3\" ||\n    t.indexOf(\"If one is not already present, create a new connection with your root user.\") === 0 ||\n    t === \"Step 4\"\n  ) {\n    toDelete.push(paragraphs.items[i]);\n  }\n}\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n\n// --- 3. Update the (now promoted) SQL-script paragraph -------------------\nconst openRange = body.search(\"Once this is opened, open a new SQL script file.\", { matchCase: true });\nopenRange.load(\"items\");\nawait context.sync();\nif (openRange.items.length > 0) {\n  openRange.items[0].insertText(\n    \"Once this is opened, open (not create!) a new SQL script file.\",\n    \"Replace\"\n  );\n}\n\nconst fileNamedRange = body.search(\"AirBooksUser.sql\u201d. Once this is opened in the workbench\", { matchCase: true });\nfileNamedRange.load(\"items\");\nawait context.sync();\nif (fileNamedRange.items.length > 0) {\n  fileNamedRange.items[0].insertText(\n    \"AirBooksUser.sql\u201d, which is a file located in the same directory as this document. Once this is opened in the workbench\",\n    \"Replace\"\n  );\n}\n\n// --- 4. Renumber the trailing heading: \"Step 5\" -> \"Step 3\" ---------------\nconst step5 = body.search(\"Step 5\", { matchCase: true });\nstep5.load(\"items\");\nawait context.sync();\nif (step5.items.length > 0) {\n  step5.items[0].insertText(\"Step 3\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the documented edit to READ_THIS_FIRST.docx:\n#  1. \"During the install...\" paragraph: drop \"Client only or \" and append a\n#     new sentence about creating a root account/password.\n#  2. Remove the old \"Step 2\" body paragraph (the MySQL-Client-setup\n#     hyperlink blurb) together with the whole \"Step 3\" / \"Step 4\" sections,\n#     because their replacement content already exists lower in the\n#     document (the paragraphs that used to live under \"Step 4\") - removing\n#     the now-obsolete paragraphs promotes that content up to sit right\n#     under \"Step 2\".\n#  3. Tweak the (promoted) \"Once this is opened, open a new SQL script\n#     file...\" paragraph: insert \"(not create!)\" and the \"located in the\n#     same directory...\" clause.\n#  4. Renumber the final remaining Heading2 (\"Step 5\") to \"Step 3\".\n\n$d = $word.ActiveDocument\n\n# wdReplace enum values used below with Find.Execute's positional ReplaceWith/Replace args.\n$wdReplaceNone = 0\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n# --- 1. \"During the install...\" paragraph --------------------------------\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Execute(\"Client only or \", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"\", $wdReplaceOne)\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Execute(\"Full must be selected. \", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"Full must be selected. Make sure that in the setup process you create a root account and password for the account.\", $wdReplaceOne)\n\n# --- 2. Drop the obsolete paragraphs/headings -----------------------------\n# Right after the \"Step 2\" heading (paragraph 5) the current document has:\n#   6: \"Once MySQL has been installed, follow the direction at this website...\"\n#   7: \"Step 3\"\n#   8: \"If one is not already present, create a new connection...\"\n#   9: \"Step 4\"\n# All four are removed (highest index first so earlier indices stay valid);\n# the paragraphs that used to sit beneath \"Step 4\" shift up to immediately\n# follow \"Step 2\".\n$d.Paragraphs.Item(9).Range.Delete()\n$d.Paragraphs.Item(8).Range.Delete()\n$d.Paragraphs.Item(7).Range.Delete()\n$d.Paragraphs.Item(6).Range.Delete()\n\n# --- 3. Update the (now promoted) SQL-script paragraph --------------------\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Execute(\"open a new SQL script file\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"open (not create!) a new SQL script file\", $wdReplaceOne)\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Execute(\"AirBooksUser.sql\u201d. Once this is opened in the workbench\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"AirBooksUser.sql\u201d, which is a file located in the same directory as this document. Once this is opened in the workbench\", $wdReplaceOne)\n\n# --- 4. Renumber the trailing heading: \"Step 5\" -> \"Step 3\" ---------------\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Execute(\"Step 5\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"Step 3\", $wdReplaceOne)\n"}
